$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.589.17"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.896.71"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.691"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.25"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "57.41"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +10.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.357"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0755"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("E12").Value = "  +1.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.55"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +10.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.803"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +11.17%  "

$ws.Range("D15").Value = "2.172.58"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("E16").Value = "  +2.55%  "

$ws.Range("D17").Value = "1.889.01"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "35.598.97"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").Value = "0.0₃0833"
$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "246.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.66"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.74%  "

$ws.Range("E28").Value = "  +2.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.36"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("E31").Value = "  +3.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0602"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.98%  "

$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +14.93%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.48"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -16.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.855"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0730"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.46%  "

$ws.Range("E40").Value = "  +6.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.39"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.13"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.38%  "

$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.88"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +14.88%  "

$ws.Range("D45").Value = "1.322.51"
$ws.Range("E45").Value = "  +2.87%  "

$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0811"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.27%  "

$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.73"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.75%  "
